$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison": update MyForecast (column D) values ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$wsForecast.Range("D9").Value = 33
$wsForecast.Range("D10").Value = 32
$wsForecast.Range("D13").Value = 24
$wsForecast.Range("D14").Value = 23
$wsForecast.Range("D15").Value = 30

# --- Sheet "Summary": update derived totals/min (column B), stored as text ---
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "518"
$wsSummary.Range("B9").ClearFormats()

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "286"
$wsSummary.Range("B10").ClearFormats()

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "23"
$wsSummary.Range("B14").ClearFormats()
